# v0.1 add + pagedown
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two lowest rows of the table (former rows 4 & 5:
# 移动端后台管理3 / 移动端后台管理4), shrinking the table to 3 rows.
$ws.Rows("4:5").Delete()

# Update the weights first, then the work-item names, so new shared
# strings are appended in the same order the source workbook used
# (0.3, 0.2, 后台系统优化, 财务报表优化二期).
$ws.Range("D2").Value = "0.3"
$ws.Range("D3").Value = "0.2"

$ws.Range("B2").Value = "后台系统优化"
$ws.Range("B3").Value = "财务报表优化二期"

# Move the active selection down, past the now-empty rows (page-down).
$ws.Range("A4:XFD7").Select()
